# Generate Report for Handback
# Row 8 (cd8c76c6-033a-410e-8a69-51b7e1d7f1b3) transitions from "Ready for
# handoff" to a failed handback transform, with an explanatory error detail
# recorded for both the zh-cn and de-de locales.

$wb = $excel.ActiveWorkbook

$statusText = "Handback transform failed"

$zhError = "Handback file name: 13dmwtmh.mek is different with handoff file name: cd8c76c6-033a-410e-8a69-51b7e1d7f1b3.0412de8136ff298f912f50140510bff0f4b3ef64.zh-cn."
$deError = "Handback file name: 13dmwtmh.mek is different with handoff file name: cd8c76c6-033a-410e-8a69-51b7e1d7f1b3.0412de8136ff298f912f50140510bff0f4b3ef64.de-de."

# --- Overview sheet: row 8 Status column shared by both locales ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E8").Value = $statusText
$wsOverview.Range("F8").Value = $statusText

# --- zh-cn sheet: row 8 Status + Error Detail ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C8").Value = $statusText
$wsZh.Range("R8").Value = $zhError

# --- de-de sheet: row 8 Status + Error Detail ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C8").Value = $statusText
$wsDe.Range("R8").Value = $deError
